$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -2
    3  = 1
    4  = -1
    5  = -3
    6  = 14
    7  = 13
    9  = -1
    10 = -2
    11 = -1
    12 = 1
    13 = 2
    14 = 5
    15 = 10
    16 = 4
    17 = -1
    18 = 2
    19 = 3
    20 = 2
    21 = 3
    22 = -1
    24 = -1
    25 = -7
    26 = 3
    27 = -3
    28 = -1
    29 = 3
    30 = 1
    31 = -2
    32 = 3
    33 = 4
    34 = 1
    35 = -1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
